$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price (column D) and 1h volume % change (column E) refresh.
# Column D values are stored as text (e.g. "43.199.77" uses dots as thousands
# separators, not a valid number), so we force Text format before assigning to
# avoid Excel auto-converting them to numbers / dropping trailing zeros, then
# restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.199.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.58%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.308.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.511'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.05%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.513'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.70'
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '

$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.96'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.664.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.307.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.800'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.063.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("E19").Value = '  +9.99%  '

$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("E21").Value = '  +1.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("E24").Value = '  +8.49%  '

$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '170.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.57%  '

$ws.Range("E30").Value = '  -0.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.18%  '

$ws.Range("E33").Value = '  +2.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.77'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.42%  '

$ws.Range("E35").Value = '  +0.55%  '

$ws.Range("E36").Value = '  -1.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0693'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.81%  '

$ws.Range("E38").Value = '  +1.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.80'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.83'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.04%  '

$ws.Range("E41").Value = '  +1.05%  '

$ws.Range("E42").Value = '  +3.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.987.22'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  -4.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.24'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.76%  '

$ws.Range("E49").Value = '  +4.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.531.15'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.28%  '
